$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing A3 date (was 41554 / 2013-10-07, should be 41553 / 2013-10-06)
$ws.Range("A3").Value = 41553

# Add the new row 4 entries
$ws.Range("A4").Value = 41554
$ws.Range("B4").Value = 0.14583333333333334

# Copy the date/time formatting from row 3 down to row 4 so the new cells
# pick up the same number formats/styles as the rest of the table.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)   # xlPasteFormats

# Underline the "Quantidade de horas" time column (B2:B4)
$ws.Range("B2:B4").Font.Underline = 2   # xlUnderlineStyleSingle

# Update the current selection to match the widened table
$ws.Range("B2:B4").Select()
